$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new column E
$ws.Range("E1").Value = "points"

# Fill column E with cycling 1..5 values for rows 2..104
for ($row = 2; $row -le 104; $row++) {
    $val = (($row - 2) % 5) + 1
    $ws.Cells.Item($row, 5).Value = $val
}
